# Situation sketch.pptx — insert a new, blank slide as the 2nd slide of the
# deck (directly after the current title/intro slide). This mirrors what the
# author did in PowerPoint: Right click slide 1 in the pane -> New Slide,
# picking the "Blank" layout, which pushes all the following slides down by
# one position.

$p = $ppt.ActivePresentation

# The "blank" layout is CustomLayout #7 ("Leeg") on the single slide master
# used by this deck - every other blank-content slide in the presentation
# (slides 1, 3, 4, 5) already points at this same layout.
$master = $p.Slides.Item(1).Master
$blankLayout = $master.CustomLayouts.Item(7)

# Insert the new slide at index 2, i.e. right after the existing first slide,
# shifting the former slides 2-5 down to positions 3-6.
$newSlide = $p.Slides.AddSlide(2, $blankLayout)
